# Price update for 2026-02-07
# Appends a new tracked-price row (Date, Price, Discount, Incredible) to
# the bottom of the price-history sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Locate the first empty row below the existing data in column A.
$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row
$newRow = $lastRow + 1

# New row values. These look like dates/numbers, but the tracker stores
# every column as plain text, so they must land as literal text rather
# than being auto-converted to a date serial / number by Excel.
$newValues = @("2026-02-07", "168000", "58", "0")

for ($i = 0; $i -lt $newValues.Length; $i++) {
    $col = $i + 1

    # Use a scratch cell far outside the used range to turn the literal
    # string into a text formula result, then copy/paste-special it as a
    # value. That lands a plain text value in the target cell without
    # Excel's "looks like a date/number" auto-detection kicking in, and
    # without leaving behind any custom number-format/style.
    $scratch = $ws.Cells.Item(500, 500 + $i)
    $scratch.Formula = '="' + $newValues[$i] + '"'
    $scratch.Copy()

    $target = $ws.Cells.Item($newRow, $col)
    $target.PasteSpecial(-4163)  # xlPasteValues

    $scratch.Clear()
}

$excel.CutCopyMode = 0
